$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# --- Row 2 updates ---
$ws.Range("D2").Value = -0.644
$ws.Range("G2").Value = -0.5384615384615385
$ws.Range("H2").Value = -0.5384615384615385
$ws.Range("I2").Value = -2.923076923076923
$ws.Range("J2").Value = -2.923076923076923
$ws.Range("K2").Value = -0.043
$ws.Range("L2").Value = -3.307692307692307
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 0.026
$ws.Range("V2").Value = 0.02113821138211382
$ws.Range("W2").Value = -0.05381727158948685
$ws.Range("X2").Value = 0.1433325504561864
$ws.Range("Y2").Value = -0.1971498220456732
$ws.Range("Z2").Value = 0.005110062893081761
$ws.Range("AA2").Value = -0.01493710691823899
$ws.Range("AB2").Value = 0.1038502895939482
$ws.Range("AC2").Value = -0.1187873965121872
$ws.Range("AG2").Value = 1.794
$ws.Range("AH2").Value = 0.5967213114754099
$ws.Range("AI2").Value = 0.7315112540192926
$ws.Range("AJ2").Value = 0.5932539682539683
$ws.Range("AK2").Value = 0.7286758732737612
$ws.Range("AL2").Value = 0.005
$ws.Range("AM2").Value = 0.005
$ws.Range("AO2").Value = -7.6
$ws.Range("AQ2").Value = -7.6

# --- Row 3 updates ---
$ws.Range("D3").Value = -0.644
$ws.Range("G3").Value = -0.5384615384615385
$ws.Range("H3").Value = -0.5384615384615385
$ws.Range("I3").Value = -2.923076923076923
$ws.Range("J3").Value = -2.923076923076923
$ws.Range("K3").Value = -0.043
$ws.Range("L3").Value = -3.307692307692307
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 0.026
$ws.Range("V3").Value = 0.02113821138211382
$ws.Range("W3").Value = -0.05381727158948685
$ws.Range("X3").Value = 0.1433325504561864
$ws.Range("Y3").Value = -0.1971498220456732
$ws.Range("Z3").Value = 0.005110062893081761
$ws.Range("AA3").Value = -0.01493710691823899
$ws.Range("AB3").Value = 0.1038502895939482
$ws.Range("AC3").Value = -0.1187873965121872
$ws.Range("AG3").Value = 1.794
$ws.Range("AH3").Value = 0.5967213114754099
$ws.Range("AI3").Value = 0.7315112540192926
$ws.Range("AJ3").Value = 0.5932539682539683
$ws.Range("AK3").Value = 0.7286758732737612
$ws.Range("AL3").Value = 0.005
$ws.Range("AM3").Value = 0.005
$ws.Range("AO3").Value = -7.6
$ws.Range("AQ3").Value = -7.6
